# Addition of the Section "Stats".
#
# The points-table header row (B1:G1) is reworded/re-cased from the terse
# lower-case abbreviations (" pld", " won", " lost", " tied", " nrr", " pts")
# to proper, trimmed headings: "Pld", "Won", "Lost", "Tied", "Net RR", "Pts".
#
# (Writing these new strings naturally appends them to the end of the shared-
# string table and leaves the team-name strings, which are unchanged, at
# their original relative order - matching the workbook's new string table.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Pld"
$ws.Range("C1").Value = "Won"
$ws.Range("D1").Value = "Lost"
$ws.Range("E1").Value = "Tied"
$ws.Range("F1").Value = "Net RR"
$ws.Range("G1").Value = "Pts"

# Select the full stats table (mirrors the workbook's
# "points_table_17" defined name, Sheet1!$A$1:$G$9) so the saved view
# reflects that range as selected.
[void]$ws.Range("A1:G9").Select()
